# Generate Report for Handback
# The 1f22a44e-...md file has now been handed back (in sync with en-US),
# swapping places with 819a2cbc-...md (which already was handed back),
# and the cb396497-...md dependency row now mirrors the 1f22a44e handoff/
# handback file details.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("A3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("B4").Value2 = "Handed back: in sync with en-US"
$ws.Range("C4").Value2 = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("C2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.zh-cn.xlf"
$ws.Range("D2").Value2 = "2016-03-02 10:08:14"
$ws.Range("E2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("F2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.zh-cn.xlf"
$ws.Range("G2").Value2 = "2016-03-02 10:09:00"

$ws.Range("A3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.2cc76cca37dd752ca70b144f95ca1920398c31f5.zh-cn.xlf"
$ws.Range("D3").Value2 = "2016-03-02 10:06:23"
$ws.Range("E3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.md"
$ws.Range("F3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.2cc76cca37dd752ca70b144f95ca1920398c31f5.zh-cn.xlf"
$ws.Range("G3").Value2 = "2016-03-02 10:07:10"

$ws.Range("B4").Value2 = "Handed back: in sync with en-US"
$ws.Range("E4").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("F4").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.zh-cn.xlf"
$ws.Range("G4").Value2 = "2016-03-02 10:09:00"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("C2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.de-de.xlf"
$ws.Range("D2").Value2 = "2016-03-02 10:08:25"
$ws.Range("E2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("F2").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.de-de.xlf"
$ws.Range("G2").Value2 = "2016-03-02 10:09:19"

$ws.Range("A3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.2cc76cca37dd752ca70b144f95ca1920398c31f5.de-de.xlf"
$ws.Range("D3").Value2 = "2016-03-02 10:06:35"
$ws.Range("E3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.md"
$ws.Range("F3").Value2 = "819a2cbc-9f1b-4372-9acc-939d459cfb78.2cc76cca37dd752ca70b144f95ca1920398c31f5.de-de.xlf"
$ws.Range("G3").Value2 = "2016-03-02 10:07:31"

$ws.Range("B4").Value2 = "Handed back: in sync with en-US"
$ws.Range("E4").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.md"
$ws.Range("F4").Value2 = "1f22a44e-c531-4882-a642-8fa27d3358fd.535d034ccf8fc5b3c51711cb815532455c9a7486.de-de.xlf"
$ws.Range("G4").Value2 = "2016-03-02 10:09:19"
